# Trade #16 closed at 2026-02-16 22:58:55 - base_strategy UP +0.000%
#
# Appends the new trade-log row (row 17) to both the "All Trades" sheet and
# the per-strategy "base_strategy" sheet. Both sheets keep an identical
# running log, so the same row is written to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Append after the last existing trade-log row (row 16 -> new row 17).
    $row = $ws.UsedRange.Rows.Count() + 1

    $ws.Range("A$row").Value = 16

    # The Date column holds plain "YYYY-MM-DD" text in this log, not a real
    # date. Prefix with an apostrophe so Excel stores it as text (matching
    # every other row already in the sheet) instead of silently
    # reinterpreting it as a date serial value.
    $ws.Range("B$row").Value = "'2026-02-16"
    $ws.Range("C$row").Value = "22:58:55"

    $ws.Range("D$row").Value = "base_strategy"
    $ws.Range("E$row").Value = "UP"
    $ws.Range("F$row").Value = 0.5
    $ws.Range("G$row").Value = ""
    $ws.Range("H$row").Value = "OPEN"
    $ws.Range("I$row").Value = 0
    $ws.Range("J$row").Value = 0
    $ws.Range("K$row").Value = 100
    $ws.Range("L$row").Value = 0
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0.6
    $ws.Range("O$row").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P$row").Value = ""
    $ws.Range("Q$row").Value = 0
}
